# Add separability results in csv
# Round the ConvexHullArea values (column D, rows 2-118) to the nearest
# whole number, replacing the full-precision floating point values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 118 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value()
    if ($val -ne $null) {
        $cell.Value = [math]::Round($val, 0)
    }
}
